# Session 5: Dynamic Programming
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F4: grade for Session 5 (Dynamic Prog.)
$ws.Range("F4").Value = 6.5

# F5: feedback comment for Session 5 (Dynamic Prog.)
$ws.Range("F5").Value = "Dynamic table was not created. Stack memory is around O(1) for dynamic programming and O(n) for the recursive version of the algorithm (taking into account the height of the tree of states). Complexy for recursive scheme is O(3^n) since we need to always do the 3 calls to guarantee a correct solution. The recursive implemention fails with stack overflow (too many recursive calls) and the dynamic programming fails because we are consuming too much memory to allocate the table"

# Update the active selection to match the author's edit
$ws.Range("G5:G12").Select()
